$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update FechaSiniestro (G2) and NroPoliza (E2) with the new inspection values.
# Leading apostrophe forces text (quote-prefixed) entry, matching the
# existing cell formatting used for these columns.
$ws.Range("G2").Value = "'19/05/2021"
$ws.Range("E2").Value = "'11111003252"

# Update the selected cell to reflect the new active selection
$ws.Range("E6").Select()
